$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values regenerated for the "filter save games" changes (rows 2-9, cols B-G)
$data = @{
    2 = @(3.230985683306322, 1.667794583268128, 3.900430680208489, 0.496779210170732, 1, 9.295990156953671)
    3 = @(1.459612070389937, 1.667794583268128, 0.1575252929769615, 0.496779210170732, 0, 3.781711156805759)
    4 = @(3.230985683306322, 1.667794583268128, 3.900430680208489, 0.496779210170732, 1, 9.295990156953671)
    5 = @(0.6753301551942219, 1.667794583268128, 0.8054896365839992, 8.660232485948974, 0, 11.80884686099532)
    6 = @(3.230985683306322, 1.667794583268128, 3.900430680208489, 8.660232485948974, 0, 17.45944343273191)
    7 = @(0.6753301551942219, 0.3127903958511391, 3.900430680208489, 645.3272768299601, 1, 650.2158280612139)
    8 = @(0.003994804209775715, 0.3127903958511391, 3.900430680208489, 645.3272768299601, 0, 649.5444927102294)
    9 = @(3.230985683306322, 1.667794583268128, 3.900430680208489, 0.496779210170732, 1, 9.295990156953671)
}

foreach ($r in $data.Keys) {
    $vals = $data[$r]
    $ws.Cells.Item($r, 2).Value = $vals[0]
    $ws.Cells.Item($r, 3).Value = $vals[1]
    $ws.Cells.Item($r, 4).Value = $vals[2]
    $ws.Cells.Item($r, 5).Value = $vals[3]
    $ws.Cells.Item($r, 6).Value = $vals[4]
    $ws.Cells.Item($r, 7).Value = $vals[5]
}
